$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12 data: date 1.4.24 entry for Kiosk / Eintritte
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A12").Value = 104.24
$ws.Range("B12").Value = "Spez 1"
$ws.Range("C12").Value = "Leibniz Kekse Dschungel"
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 2

$ws.Range("A13").Select()
